$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.393.66"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "3.487.61"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "609.38"
$ws.Range("E5").Value = "  +4.67%  "
$ws.Range("D6").Value = "185.93"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -3.57%  "
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("D12").Value = "0.0000308"
$ws.Range("E12").Value = "  -3.61%  "
$ws.Range("D13").Value = "9.52"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "4.037.90"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "600.37"
$ws.Range("E15").Value = "  +4.46%  "
$ws.Range("D16").Value = "69.458.53"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").Value = "18.81"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "12.52"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").Value = "3.483.22"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "0.984"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "17.09"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").Value = "105.96"
$ws.Range("E23").Value = "  +12.34%  "
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("E25").Value = "  +3.37%  "
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("D28").Value = "9.68"
$ws.Range("E28").Value = "  +4.72%  "
$ws.Range("D29").Value = "33.59"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").Value = "4.09"
$ws.Range("E31").Value = "  +17.31%  "
$ws.Range("D32").Value = "12.40"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").Value = "63.25"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "3.17"
$ws.Range("E35").Value = "  -6.78%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "521.49"
$ws.Range("E37").Value = "  -5.25%  "
$ws.Range("E38").Value = "  +7.12%  "
$ws.Range("D39").Value = "3.612.78"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("D41").Value = "36.69"
$ws.Range("E41").Value = "  -3.42%  "
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "0.0460"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("D48").Value = "8.80"
$ws.Range("E48").Value = "  -5.77%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "1.35"
$ws.Range("E50").Value = "  -10.55%  "
$ws.Range("E51").Value = "  -8.65%  "
